$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_WVR = $wb.Worksheets.Item("WVR")

# ALC row 97 (G97=19885)
$ws_ALC.Range("H97").Value = 398.33334
$ws_ALC.Range("I97").Value = 0
$ws_ALC.Range("J97").Value = 398.33334
$ws_ALC.Range("K97").Value = 0
$ws_ALC.Range("L97").Value = 1195.00002
$ws_ALC.Range("M97").Value = ""
$ws_ALC.Range("N97").Value = -2187.00002

# ALC row 105 (G105=18668)
$ws_ALC.Range("H105").Value = 40000
$ws_ALC.Range("I105").Value = 0
$ws_ALC.Range("J105").Value = 40000
$ws_ALC.Range("K105").Value = 0
$ws_ALC.Range("L105").Value = 40000
$ws_ALC.Range("N105").Value = -46988

# ALC row 121 (G121=39731)
$ws_ALC.Range("H121").Value = 1912.4546
$ws_ALC.Range("I121").Value = 0
$ws_ALC.Range("J121").Value = 1912.4546
$ws_ALC.Range("K121").Value = 0
$ws_ALC.Range("L121").Value = 5737.3638
$ws_ALC.Range("N121").Value = -9231.363799999999

# ARM row 45 (G45=27714)
$ws_ARM.Range("H45").Value = 2349.1428
$ws_ARM.Range("I45").Value = 2326.3635
$ws_ARM.Range("J45").Value = 2432.6667
$ws_ARM.Range("K45").Value = 2326.3635
$ws_ARM.Range("L45").Value = 2432.6667
$ws_ARM.Range("M45").Value = -1949.3635
$ws_ARM.Range("N45").Value = -3186.6667

# ARM row 61 (G61=43999)
$ws_ARM.Range("H61").Value = 6309.8887
$ws_ARM.Range("I61").Value = 6309.8887
$ws_ARM.Range("J61").Value = 0
$ws_ARM.Range("K61").Value = 6309.8887
$ws_ARM.Range("L61").Value = 0
$ws_ARM.Range("M61").Value = -6097.8887

# ARM row 74 (G74=44000)
$ws_ARM.Range("H74").Value = 2219.923
$ws_ARM.Range("I74").Value = 2460.2727
$ws_ARM.Range("J74").Value = 898
$ws_ARM.Range("K74").Value = 2460.2727
$ws_ARM.Range("L74").Value = 898
$ws_ARM.Range("M74").Value = -1586.2727

# ARM row 77 (G77=44000)
$ws_ARM.Range("H77").Value = 2219.923
$ws_ARM.Range("I77").Value = 2460.2727
$ws_ARM.Range("J77").Value = 898
$ws_ARM.Range("K77").Value = 12301.3635
$ws_ARM.Range("L77").Value = 4490
$ws_ARM.Range("M77").Value = -7933.363499999999

# ARM row 96 (G96=18207)
$ws_ARM.Range("H96").Value = 40000
$ws_ARM.Range("I96").Value = 0
$ws_ARM.Range("J96").Value = 40000
$ws_ARM.Range("K96").Value = 0
$ws_ARM.Range("L96").Value = 40000
$ws_ARM.Range("N96").Value = -45492

# ARM row 102 (G102=19945)
$ws_ARM.Range("H102").Value = 2749.25
$ws_ARM.Range("I102").Value = 2749.25
$ws_ARM.Range("J102").Value = 0
$ws_ARM.Range("K102").Value = 2749.25
$ws_ARM.Range("L102").Value = 0
$ws_ARM.Range("M102").Value = -1127.25

# ARM row 132 (G132=43997)
$ws_ARM.Range("H132").Value = 3147.4443
$ws_ARM.Range("I132").Value = 2596.7693
$ws_ARM.Range("J132").Value = 4579.2
$ws_ARM.Range("K132").Value = 7790.3079
$ws_ARM.Range("L132").Value = 13737.6
$ws_ARM.Range("M132").Value = -5260.3079

# ARM row 136 (G136=43999)
$ws_ARM.Range("H136").Value = 6309.8887
$ws_ARM.Range("I136").Value = 6309.8887
$ws_ARM.Range("J136").Value = 0
$ws_ARM.Range("K136").Value = 18929.6661
$ws_ARM.Range("L136").Value = 0
$ws_ARM.Range("M136").Value = -16379.6661

# BSM row 94 (G94=19939)
$ws_BSM.Range("H94").Value = 401
$ws_BSM.Range("I94").Value = 401
$ws_BSM.Range("J94").Value = 0
$ws_BSM.Range("K94").Value = 401
$ws_BSM.Range("L94").Value = 0
$ws_BSM.Range("M94").Value = 50

# BSM row 105 (G105=19947)
$ws_BSM.Range("H105").Value = 2071
$ws_BSM.Range("I105").Value = 3144
$ws_BSM.Range("J105").Value = 998
$ws_BSM.Range("K105").Value = 3144
$ws_BSM.Range("L105").Value = 998
$ws_BSM.Range("M105").Value = -1397

# CRP row 31 (G31=44023)
$ws_CRP.Range("H31").Value = 932.94116
$ws_CRP.Range("I31").Value = 983
$ws_CRP.Range("J31").Value = 897.9
$ws_CRP.Range("K31").Value = 983
$ws_CRP.Range("L31").Value = 897.9
$ws_CRP.Range("M31").Value = -688
$ws_CRP.Range("N31").Value = -1487.9

# CRP row 34 (G34=44023)
$ws_CRP.Range("H34").Value = 932.94116
$ws_CRP.Range("I34").Value = 983
$ws_CRP.Range("J34").Value = 897.9
$ws_CRP.Range("K34").Value = 983
$ws_CRP.Range("L34").Value = 897.9
$ws_CRP.Range("M34").Value = -781
$ws_CRP.Range("N34").Value = -1301.9

# CRP row 64 (G64=10610)
$ws_CRP.Range("H64").Value = 69271
$ws_CRP.Range("I64").Value = 0
$ws_CRP.Range("J64").Value = 69271
$ws_CRP.Range("K64").Value = 0
$ws_CRP.Range("L64").Value = 69271
$ws_CRP.Range("N64").Value = -69767

# CRP row 67 (G67=10610)
$ws_CRP.Range("H67").Value = 69271
$ws_CRP.Range("I67").Value = 0
$ws_CRP.Range("J67").Value = 69271
$ws_CRP.Range("K67").Value = 0
$ws_CRP.Range("L67").Value = 69271
$ws_CRP.Range("N67").Value = -70987

# CRP row 99 (G99=36198)
$ws_CRP.Range("H99").Value = 0
$ws_CRP.Range("I99").Value = 0
$ws_CRP.Range("J99").Value = 0
$ws_CRP.Range("K99").Value = 0
$ws_CRP.Range("L99").Value = 0
$ws_CRP.Range("N99").Value = ""

# CRP row 105 (G105=19928)
$ws_CRP.Range("H105").Value = 3299.5
$ws_CRP.Range("I105").Value = 3266.3333
$ws_CRP.Range("J105").Value = 3332.6667
$ws_CRP.Range("K105").Value = 3266.3333
$ws_CRP.Range("L105").Value = 3332.6667
$ws_CRP.Range("M105").Value = -1519.3333
$ws_CRP.Range("N105").Value = -6826.6667

# CRP row 106 (G106=18661)
$ws_CRP.Range("H106").Value = 100000
$ws_CRP.Range("I106").Value = 0
$ws_CRP.Range("J106").Value = 100000
$ws_CRP.Range("K106").Value = 0
$ws_CRP.Range("L106").Value = 100000
$ws_CRP.Range("N106").Value = -102524

# CRP row 126 (G126=36198)
$ws_CRP.Range("H126").Value = 0
$ws_CRP.Range("I126").Value = 0
$ws_CRP.Range("J126").Value = 0
$ws_CRP.Range("K126").Value = 0
$ws_CRP.Range("L126").Value = 0
$ws_CRP.Range("N126").Value = ""

# CRP row 132 (G132=44019)
$ws_CRP.Range("H132").Value = 4386.25
$ws_CRP.Range("I132").Value = 2774.25
$ws_CRP.Range("J132").Value = 5998.25
$ws_CRP.Range("K132").Value = 8322.75
$ws_CRP.Range("L132").Value = 17994.75
$ws_CRP.Range("M132").Value = -5792.75

# CUL row 5 (G5=43974)
$ws_CUL.Range("H5").Value = 1000
$ws_CUL.Range("I5").Value = 0
$ws_CUL.Range("J5").Value = 1000
$ws_CUL.Range("K5").Value = 0
$ws_CUL.Range("L5").Value = 3000
$ws_CUL.Range("M5").Value = ""
$ws_CUL.Range("N5").Value = -3224

# CUL row 22 (G22=4697)
$ws_CUL.Range("H22").Value = 1250
$ws_CUL.Range("I22").Value = 0
$ws_CUL.Range("J22").Value = 1250
$ws_CUL.Range("K22").Value = 0
$ws_CUL.Range("L22").Value = 3750
$ws_CUL.Range("N22").Value = -4088

# CUL row 24 (G24=4690)
$ws_CUL.Range("H24").Value = 2001
$ws_CUL.Range("I24").Value = 0
$ws_CUL.Range("J24").Value = 2001
$ws_CUL.Range("K24").Value = 0
$ws_CUL.Range("L24").Value = 6003
$ws_CUL.Range("N24").Value = -6463

# CUL row 25 (G25=4709)
$ws_CUL.Range("H25").Value = 500
$ws_CUL.Range("I25").Value = 500
$ws_CUL.Range("J25").Value = 0
$ws_CUL.Range("K25").Value = 1500
$ws_CUL.Range("L25").Value = 0
$ws_CUL.Range("M25").Value = -1331

# CUL row 27 (G27=4697)
$ws_CUL.Range("H27").Value = 1250
$ws_CUL.Range("I27").Value = 0
$ws_CUL.Range("J27").Value = 1250
$ws_CUL.Range("K27").Value = 0
$ws_CUL.Range("L27").Value = 3750
$ws_CUL.Range("N27").Value = -3954

# CUL row 30 (G30=4709)
$ws_CUL.Range("H30").Value = 500
$ws_CUL.Range("I30").Value = 500
$ws_CUL.Range("J30").Value = 0
$ws_CUL.Range("K30").Value = 1500
$ws_CUL.Range("L30").Value = 0
$ws_CUL.Range("M30").Value = -1398

# CUL row 59 (G59=4694)
$ws_CUL.Range("H59").Value = 25000
$ws_CUL.Range("I59").Value = 0
$ws_CUL.Range("J59").Value = 25000
$ws_CUL.Range("K59").Value = 0
$ws_CUL.Range("L59").Value = 75000
$ws_CUL.Range("M59").Value = ""
$ws_CUL.Range("N59").Value = -76080

# CUL row 74 (G74=12859)
$ws_CUL.Range("H74").Value = 7500
$ws_CUL.Range("I74").Value = 0
$ws_CUL.Range("J74").Value = 7500
$ws_CUL.Range("K74").Value = 0
$ws_CUL.Range("L74").Value = 22500
$ws_CUL.Range("N74").Value = -24622

# CUL row 77 (G77=12859)
$ws_CUL.Range("H77").Value = 7500
$ws_CUL.Range("I77").Value = 0
$ws_CUL.Range("J77").Value = 7500
$ws_CUL.Range("K77").Value = 0
$ws_CUL.Range("L77").Value = 67500
$ws_CUL.Range("N77").Value = -78108

# CUL row 124 (G124=36040)
$ws_CUL.Range("H124").Value = 5000
$ws_CUL.Range("I124").Value = 0
$ws_CUL.Range("J124").Value = 5000
$ws_CUL.Range("K124").Value = 0
$ws_CUL.Range("L124").Value = 15000
$ws_CUL.Range("N124").Value = -24820

# CUL row 126 (G126=36045)
$ws_CUL.Range("H126").Value = 7500
$ws_CUL.Range("I126").Value = 0
$ws_CUL.Range("J126").Value = 7500
$ws_CUL.Range("K126").Value = 0
$ws_CUL.Range("L126").Value = 22500
$ws_CUL.Range("N126").Value = -32380

# CUL row 132 (G132=43972)
$ws_CUL.Range("H132").Value = 1833.3334
$ws_CUL.Range("I132").Value = 1833.3334
$ws_CUL.Range("J132").Value = 0
$ws_CUL.Range("K132").Value = 16500.0006
$ws_CUL.Range("L132").Value = 0
$ws_CUL.Range("M132").Value = -13970.0006

# CUL row 133 (G133=44073)
$ws_CUL.Range("H133").Value = 10000
$ws_CUL.Range("I133").Value = 10000
$ws_CUL.Range("J133").Value = 0
$ws_CUL.Range("K133").Value = 30000
$ws_CUL.Range("L133").Value = 0
$ws_CUL.Range("M133").Value = -24940
$ws_CUL.Range("N133").Value = ""

# CUL row 134 (G134=44074)
$ws_CUL.Range("H134").Value = 1833.3334
$ws_CUL.Range("I134").Value = 1833.3334
$ws_CUL.Range("J134").Value = 0
$ws_CUL.Range("K134").Value = 5500.0002
$ws_CUL.Range("L134").Value = 0
$ws_CUL.Range("M134").Value = -430.0002000000004

# CUL row 135 (G135=43974)
$ws_CUL.Range("H135").Value = 1000
$ws_CUL.Range("I135").Value = 0
$ws_CUL.Range("J135").Value = 1000
$ws_CUL.Range("K135").Value = 0
$ws_CUL.Range("L135").Value = 9000
$ws_CUL.Range("M135").Value = ""
$ws_CUL.Range("N135").Value = -14070

# GSM row 70 (G70=14146)
$ws_GSM.Range("H70").Value = 5500
$ws_GSM.Range("I70").Value = 5500
$ws_GSM.Range("J70").Value = 0
$ws_GSM.Range("K70").Value = 5500
$ws_GSM.Range("L70").Value = 0
$ws_GSM.Range("M70").Value = -5230

# GSM row 73 (G73=14146)
$ws_GSM.Range("H73").Value = 5500
$ws_GSM.Range("I73").Value = 5500
$ws_GSM.Range("J73").Value = 0
$ws_GSM.Range("K73").Value = 5500
$ws_GSM.Range("L73").Value = 0
$ws_GSM.Range("M73").Value = -4564

# GSM row 132 (G132=44008)
$ws_GSM.Range("H132").Value = 3525.4
$ws_GSM.Range("I132").Value = 2894.5
$ws_GSM.Range("J132").Value = 4471.75
$ws_GSM.Range("K132").Value = 8683.5
$ws_GSM.Range("L132").Value = 13415.25
$ws_GSM.Range("M132").Value = -6153.5
$ws_GSM.Range("N132").Value = -18475.25

# LTW row 46 (G46=5282)
$ws_LTW.Range("H46").Value = 2000
$ws_LTW.Range("I46").Value = 1000
$ws_LTW.Range("J46").Value = 3000
$ws_LTW.Range("K46").Value = 1000
$ws_LTW.Range("L46").Value = 3000
$ws_LTW.Range("M46").Value = -812
$ws_LTW.Range("N46").Value = -3376

# LTW row 136 (G136=44060)
$ws_LTW.Range("H136").Value = 2858.4167
$ws_LTW.Range("I136").Value = 2858.4167
$ws_LTW.Range("J136").Value = 0
$ws_LTW.Range("K136").Value = 8575.250100000001
$ws_LTW.Range("L136").Value = 0
$ws_LTW.Range("M136").Value = -6025.250100000001

# WVR row 81 (G81=12596)
$ws_WVR.Range("H81").Value = 21334.666
$ws_WVR.Range("I81").Value = 0
$ws_WVR.Range("J81").Value = 21334.666
$ws_WVR.Range("K81").Value = 0
$ws_WVR.Range("L81").Value = 42669.332
$ws_WVR.Range("N81").Value = -44791.332

# WVR row 84 (G84=12596)
$ws_WVR.Range("H84").Value = 21334.666
$ws_WVR.Range("I84").Value = 0
$ws_WVR.Range("J84").Value = 21334.666
$ws_WVR.Range("K84").Value = 0
$ws_WVR.Range("L84").Value = 213346.66
$ws_WVR.Range("N84").Value = -223954.66

# WVR row 96 (G96=19977)
$ws_WVR.Range("H96").Value = 0
$ws_WVR.Range("I96").Value = 0
$ws_WVR.Range("J96").Value = 0
$ws_WVR.Range("K96").Value = 0
$ws_WVR.Range("L96").Value = 0
$ws_WVR.Range("M96").Value = ""

# WVR row 136 (G136=44031)
$ws_WVR.Range("H136").Value = 500
$ws_WVR.Range("I136").Value = 500
$ws_WVR.Range("J136").Value = 0
$ws_WVR.Range("K136").Value = 1500
$ws_WVR.Range("L136").Value = 0
$ws_WVR.Range("M136").Value = 1050
